$wb = $excel.ActiveWorkbook

# NOTE: this COM bridge resolves worksheet variables by tab *position*, not
# stable object identity, so any variable captured before a structural
# change (Add/Move/Delete) can silently end up pointing at a different
# sheet afterwards. To stay safe, sheets are re-fetched by name immediately
# before every operation that follows a structural change.

# --- 1. Rename the three existing sheets (pure renames; no tab shifting) ---
$wb.Worksheets.Item("Example1").Name = "Матрица привлекательности"
$wb.Worksheets.Item("Example2").Name = "Матрица распределения объемов к"
$wb.Worksheets.Item("Example3").Name = "Матрица интенсивностей на часов"

# --- 2. Replace the "attractiveness" sheet's data with the new 5x5 matrix ---
$wsAttract = $wb.Worksheets.Item("Матрица привлекательности")
$wsAttract.Cells.Clear()

$attractData = New-Object 'object[,]' 5,5
$attractValues = @(
    @(0,     0.171, 0.253, 0.332, 0.244),
    @(0.121, 0,     0.268, 0.352, 0.259),
    @(0.133, 0.199, 0,     0.385, 0.283),
    @(0.146, 0.219, 0.323, 0,     0.312),
    @(0.131, 0.197, 0.291, 0.381, 0)
)
for ($r = 0; $r -lt 5; $r++) {
    for ($c = 0; $c -lt 5; $c++) {
        $attractData[$r, $c] = $attractValues[$r][$c]
    }
}
$wsAttract.Range("A1:E5").Value = $attractData
[void]$wsAttract.Range("A1").Select()

# --- 3. Empty out the other two renamed sheets, leaving their old A1:N14 selection ---
$wsVolume = $wb.Worksheets.Item("Матрица распределения объемов к")
$wsVolume.Cells.Clear()
[void]$wsVolume.Range("A1:N14").Select()

$wsIntensity = $wb.Worksheets.Item("Матрица интенсивностей на часов")
$wsIntensity.Cells.Clear()
[void]$wsIntensity.Range("A1:N14").Select()

# --- 4. Insert a brand-new blank sheet and place it first in the tab order ---
$volumeAnchor = $wb.Worksheets.Item("Матрица распределения объемов к")
$newFirst = $wb.Worksheets.Add($volumeAnchor, $null)
$newFirst.Name = "Лист1"

# --- 5. Move the "attractiveness" sheet to be last in the tab order ---
$moveMe = $wb.Worksheets.Item("Матрица привлекательности")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$moveMe.Move($null, $lastSheet)

# --- 6. Make the "volume distribution" sheet (2nd tab) the active one ---
$wb.Worksheets.Item("Матрица распределения объемов к").Activate()
[void]$wb.ActiveSheet.Range("A1:N14").Select()
